$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Parson Russell Terrier'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $casesQuery

# Reflect the author's final cursor position on the "Cases" query cell.
$ws.Range("B2").Select() | Out-Null
